# sort and export dashboard
# Insert two new participant rows into the dashboard table at the
# positions they belong (keeping the rest of the rows in their
# existing order), matching the source "before -> after" edit:
#   1. "Dave Pratikkumar Jayeshbhai" inserted right after
#      "Induni lakshika abeysinghe" (original row 18).
#   2. "Shah Pooja Sanjaybhai" inserted right after
#      "PATEL PRATIK" (original row 32, i.e. after the first
#      insertion has shifted everything down by one -> row 33).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-RowValues($Row, $Values) {
    for ($i = 0; $i -lt $Values.Length; $i++) {
        $ws.Cells.Item($Row, $i + 1).Value = $Values[$i]
    }
}

# --- Insert row 1: "Dave Pratikkumar Jayeshbhai" --------------------------
# Goes right after "Induni lakshika abeysinghe" which is row 18.
$ws.Rows.Item(18).Insert()
Set-RowValues 18 @(
    "Dave Pratikkumar Jayeshbhai",
    "9978921363",
    "140373109023",
    "ELECTRICAL",
    "8",
    "pratik.dave90@ymail.com",
    "PARUL INSTITUTE OF ENGINEERING & TECHNOLOGY",
    "WK204529"
)

# --- Insert row 2: "Shah Pooja Sanjaybhai" ---------------------------------
# Goes right after "PATEL PRATIK". Before the first insert that record was
# row 32; after inserting a row above it (row 18), it is now row 33.
$ws.Rows.Item(33).Insert()
Set-RowValues 33 @(
    "Shah Pooja Sanjaybhai",
    "9825964113",
    "130800106042",
    "080",
    "8",
    "shahpooja277@gmail.com",
    "VADODARA INSTITUTE OF ENGINEERING",
    "WK506034"
)
